$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove all existing hyperlinks first so stale refs don't linger ---
$ws.Hyperlinks.Delete()

# --- Structural changes: drop the old NETGEAR-switch data row (old row 9) ---
# and the trailing "Power Supply Search" block (old rows 15:16, now 14:15
# after the row-9 delete shifts everything up by one).
$ws.Rows("9:9").Delete()
$ws.Rows("14:15").Delete()

# --- Row 3: HP switch's URL cell (text unchanged) now gets a hyperlink, added below ---

# --- Row 5: Power Supply -> Meanwell Power Supply, new price & URL ---
$ws.Range("F5").Value = "http://www.amazon.com/Meanwell-RSP-320-5-Power-Supply-OlympianLED/dp/B00IWC2RLS/ref=sr_1_1?ie=UTF8&qid=1423685983&sr=8-1&keywords=meanwell+5v+300w"
$ws.Range("B5").Value = "Meanwell Power  Supply (5V, 60A, 300W) "
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 59.99
$ws.Range("E5").Formula = "=D5*C5"

# --- Row 7: Ethernet Cable (5-set), shorter cables, new price & URL ---
$ws.Range("F7").Value = "http://www.amazon.com/Cable-Matters-5-Color-Snagless-Ethernet/dp/B00E5I7UAG/ref=sr_1_4?ie=UTF8&qid=1423686055&sr=8-4&keywords=3ft.+ethernet+cable"
$ws.Range("B7").Value = "Ethernet Cable (5-set)"
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 10.99

# --- Row 6: fill in the previously-empty 10/100 Ethernet Adapter line ---
$ws.Range("F6").Value = "http://www.newegg.com/Product/Product.aspx?Item=9SIA2XB12C5920&cm_re=usb_ethernet-_-9SIA2XB12C5920-_-Product"
$ws.Range("B6").Value = "10/100 Ethernet Adapter (w/ linux)"
$ws.Range("C6").Value = 32
$ws.Range("D6").Value = 14.95

# --- Row 8: NETGEAR 48-Port switch takes over this slot ---
$ws.Range("B8").Value = "NETGEAR 48-Port 10/100/1000Mbps"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 396
$ws.Range("F8").Value = "http://www.amazon.com/NETGEAR-ProSAFE-GS748T-48-Port-1000Mbps/dp/B00I5W5M12/ref=sr_1_1?s=pc&ie=UTF8&qid=1422583997&sr=1-1&keywords=48+port+gigabit+switch&pebp=1422583998267&peasin=B00I5W5M12"

# --- "cables ordered" now tracks the Ethernet-Cable row, which is row 7 ---
$ws.Range("C11").Formula = "=C7*5"

# --- Re-create hyperlinks in the same order as the target workbook ---
$ws.Hyperlinks.Add($ws.Range("F4"), "http://www.amazon.com/Kingston-Digital-microSDHC-SDC4-16GBET/dp/B00DYQYLQQ/ref=sr_1_1?ie=UTF8&qid=1422292553&sr=8-1&keywords=micro+sd+16gb")
$ws.Hyperlinks.Add($ws.Range("F7"), "http://www.amazon.com/Cable-Matters-5-Color-Snagless-Ethernet/dp/B00E5I7UAG/ref=sr_1_4?ie=UTF8&qid=1423686055&sr=8-4&keywords=3ft.+ethernet+cable")
$ws.Hyperlinks.Add($ws.Range("F5"), "http://www.amazon.com/Meanwell-RSP-320-5-Power-Supply-OlympianLED/dp/B00IWC2RLS/ref=sr_1_1?ie=UTF8&qid=1423685983&sr=8-1&keywords=meanwell+5v+300w")
$ws.Hyperlinks.Add($ws.Range("F6"), "http://www.newegg.com/Product/Product.aspx?Item=9SIA2XB12C5920&cm_re=usb_ethernet-_-9SIA2XB12C5920-_-Product")
$ws.Hyperlinks.Add($ws.Range("F3"), "http://www.amazon.com/HP-J9728A-2920-48G-Switch/dp/B00BJ42JQY")
$ws.Hyperlinks.Add($ws.Range("F8"), "http://www.amazon.com/NETGEAR-ProSAFE-GS748T-48-Port-1000Mbps/dp/B00I5W5M12/ref=sr_1_1?s=pc&ie=UTF8&qid=1422583997&sr=1-1&keywords=48+port+gigabit+switch&pebp=1422583998267&peasin=B00I5W5M12")

# `Hyperlinks.Add` stamps a brand-new (duplicate) cell style onto each
# target cell instead of reusing the existing hyperlink style already on
# column F. Re-apply the original F-column hyperlink formatting (copied
# from the untouched F2 cell) so the style table doesn't grow spuriously.
$ws.Range("F2").Copy()
$ws.Range("F3:F8").PasteSpecial(-4122)

# --- Move the active selection to match the new layout ---
$ws.Range("B16").Select()
